$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Widen the "NOMBRE DEL PROYECTO / SIGLAS DEL PROYECTO" table (2nd table)
#    Table width 9102 -> 10950 (dxa); first column 4551 -> 6399 (dxa).
#    Word stores widths in points (1 pt = 20 dxa) on the COM surface.
# ---------------------------------------------------------------------------
$projTable = $d.Tables.Item(2)
$projTable.PreferredWidth = 547.5          # 10950 dxa / 20
$projTable.Columns.Item(1).Width = 319.95  # 6399 dxa / 20

# ---------------------------------------------------------------------------
# 2) Replace the project-name cell's text and formatting.
#    "APLICACIÓN MÓVIL SAN PEDRO" ->
#    "ADMINISTRACIÓN DE LA APLICACIÓN MÓVIL DE LA PANADERÍA SAN PEDRO"
#    with new run/paragraph formatting.
# ---------------------------------------------------------------------------
$cell = $projTable.Cell(2, 1)
$cellRange = $cell.Range
$textRange = $d.Range($cellRange.Start, $cellRange.End - 1)
$textRange.Delete()

$projTable2 = $d.Tables.Item(2)
$cell2 = $projTable2.Cell(2, 1)
$cell2.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:widowControl/><w:spacing w:after="160" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-ES"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/><w:color w:val="000000"/><w:lang w:val="es-ES" w:eastAsia="es-MX"/></w:rPr><w:t>ADMINISTRACIÓN DE LA APLICACIÓN MÓVIL DE LA PANADERÍA SAN PEDRO</w:t></w:r></w:p>')

# Drop the now-empty leading paragraph that InsertXML leaves behind, so the
# cell ends up with exactly one paragraph again.
$projTable3 = $d.Tables.Item(2)
$cell3 = $projTable3.Cell(2, 1)
$leadingPara = $cell3.Range.Paragraphs.Item(1)
$leadingPara.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Merge the two adjacent runs " " + "plan de recursos humanos" (same
#    formatting) into a single run " plan de recursos humanos".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" plan de recursos humanos", $true, $false, $false, $false, $false,
                         $true, 1, $false, " plan de recursos humanos", 2)
